# Reorder the player rows (rows 6-17) so that each player's Position/Team
# values travel with them, matching the new row order from the update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Kyle Kuzma", "PF", "Milwaukee Bucks"),
    @("Zach LaVine", "SG,SF", "Sacramento Kings"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Bam Adebayo", "PF,C", "Miami Heat"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Matas Buzelis", "SF,PF", "Chicago Bulls"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers"),
    @("Joel Embiid", "C", "Philadelphia 76ers")
)

$startRow = 6
for ($i = 0; $i -lt $players.Count; $i++) {
    $row = $startRow + $i
    $rec = $players[$i]
    $ws.Range("A$row").Value = $rec[0]
    $ws.Range("B$row").Value = $rec[1]
    $ws.Range("C$row").Value = $rec[2]
}
